$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: copy formatting from row 10 down to row 11 first
$ws.Range("A10:AQ10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 2036.46
$ws.Range("C11").Value = 459.11
$ws.Range("D11").Value = 196.77
$ws.Range("E11").Formula = "'"
$ws.Range("E11").ClearFormats()
$ws.Range("F11").Value = 144.15
$ws.Range("G11").Value = 1072.09
$ws.Range("H11").Value = 56.93
$ws.Range("I11").Value = 675.1900000000001
$ws.Range("J11").Value = 108.63
$ws.Range("K11").Value = 40194.1
$ws.Range("L11").Value = 42.53
$ws.Range("M11").Value = 110.57
$ws.Range("N11").Value = 250.23
$ws.Range("O11").Value = 19.28
$ws.Range("P11").Value = 843.74
$ws.Range("Q11").Value = 122.57
$ws.Range("R11").Value = 15.72
$ws.Range("S11").Value = 228.47
$ws.Range("T11").Value = 734.79
$ws.Range("U11").Value = 3940.22
$ws.Range("V11").Value = 302.3
$ws.Range("W11").Value = 2903.71
$ws.Range("X11").Value = 358.56
$ws.Range("Y11").Value = 8679.99
$ws.Range("Z11").Value = 2091.49
$ws.Range("AA11").Value = 7.14
$ws.Range("AB11").Value = 675.8
$ws.Range("AC11").Value = 407.56
$ws.Range("AD11").Value = 64.26000000000001
$ws.Range("AE11").Value = 30.13
$ws.Range("AF11").Value = 3368.13
$ws.Range("AG11").Value = 2144.59
$ws.Range("AH11").Value = 64.04000000000001
$ws.Range("AI11").Value = 156.59
$ws.Range("AJ11").Value = 223.6
$ws.Range("AK11").Value = 836.77
$ws.Range("AL11").Value = 3062.68
$ws.Range("AM11").Value = 2100.7
$ws.Range("AN11").Value = 145.18
$ws.Range("AO11").Value = 208.12
$ws.Range("AP11").Value = 1035.98
$ws.Range("AQ11").Value = 268.89

